# Append 65 new EURUSD daily/intraday OHLCV rows (rows 4993-5057) to the
# bottom of the existing price history table on Sheet1, extending the
# used range from A1:F4992 to A1:F5057.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 65,6
$data[0,0] = 45229
$data[0,1] = 1.05626
$data[0,2] = 1.05707
$data[0,3] = 1.05551
$data[0,4] = 1.05625
$data[0,5] = 11757
$data[1,0] = 45229
$data[1,1] = 1.05625
$data[1,2] = 1.05658
$data[1,3] = 1.05583
$data[1,4] = 1.05609
$data[1,5] = 15372
$data[2,0] = 45229
$data[2,1] = 1.05609
$data[2,2] = 1.05674
$data[2,3] = 1.05471
$data[2,4] = 1.05598
$data[2,5] = 29587
$data[3,0] = 45229
$data[3,1] = 1.05598
$data[3,2] = 1.06084
$data[3,3] = 1.05597
$data[3,4] = 1.06068
$data[3,5] = 41413
$data[4,0] = 45230
$data[4,1] = 1.06068
$data[4,2] = 1.06253
$data[4,3] = 1.05908
$data[4,4] = 1.06071
$data[4,5] = 74171
$data[5,0] = 45230
$data[5,1] = 1.06071
$data[5,2] = 1.06217
$data[5,3] = 1.06001
$data[5,4] = 1.0615
$data[5,5] = 22570
$data[6,0] = 45230
$data[6,1] = 1.0615
$data[6,2] = 1.06172
$data[6,3] = 1.06093
$data[6,4] = 1.06115
$data[6,5] = 6936
$data[7,0] = 45230
$data[7,1] = 1.06115
$data[7,2] = 1.06126
$data[7,3] = 1.05906
$data[7,4] = 1.05938
$data[7,5] = 28030
$data[8,0] = 45230
$data[8,1] = 1.05938
$data[8,2] = 1.06458
$data[8,3] = 1.05926
$data[8,4] = 1.06453
$data[8,5] = 34443
$data[9,0] = 45230
$data[9,1] = 1.06453
$data[9,2] = 1.06749
$data[9,3] = 1.06211
$data[9,4] = 1.06241
$data[9,5] = 60573
$data[10,0] = 45231
$data[10,1] = 1.06241
$data[10,2] = 1.06257
$data[10,3] = 1.05574
$data[10,4] = 1.05675
$data[10,5] = 90552
$data[11,0] = 45231
$data[11,1] = 1.05675
$data[11,2] = 1.05858
$data[11,3] = 1.05657
$data[11,4] = 1.05744
$data[11,5] = 29256
$data[12,0] = 45231
$data[12,1] = 1.05744
$data[12,2] = 1.05806
$data[12,3] = 1.05685
$data[12,4] = 1.05733
$data[12,5] = 10055
$data[13,0] = 45231
$data[13,1] = 1.05733
$data[13,2] = 1.05764
$data[13,3] = 1.05658
$data[13,4] = 1.05672
$data[13,5] = 18408
$data[14,0] = 45231
$data[14,1] = 1.05672
$data[14,2] = 1.05775
$data[14,3] = 1.0557
$data[14,4] = 1.05582
$data[14,5] = 28917
$data[15,0] = 45231
$data[15,1] = 1.05582
$data[15,2] = 1.05627
$data[15,3] = 1.05369
$data[15,4] = 1.05385
$data[15,5] = 47594
$data[16,0] = 45232
$data[16,1] = 1.05385
$data[16,2] = 1.05708
$data[16,3] = 1.05215
$data[16,4] = 1.05338
$data[16,5] = 80800
$data[17,0] = 45232
$data[17,1] = 1.05338
$data[17,2] = 1.05758
$data[17,3] = 1.05165
$data[17,4] = 1.05698
$data[17,5] = 91229
$data[18,0] = 45232
$data[18,1] = 1.05698
$data[18,2] = 1.0602
$data[18,3] = 1.05645
$data[18,4] = 1.06017
$data[18,5] = 14401
$data[19,0] = 45232
$data[19,1] = 1.06017
$data[19,2] = 1.06021
$data[19,3] = 1.05932
$data[19,4] = 1.05959
$data[19,5] = 17839
$data[20,0] = 45232
$data[20,1] = 1.05959
$data[20,2] = 1.06103
$data[20,3] = 1.05911
$data[20,4] = 1.06005
$data[20,5] = 33717
$data[21,0] = 45232
$data[21,1] = 1.06005
$data[21,2] = 1.06675
$data[21,3] = 1.0599
$data[21,4] = 1.06641
$data[21,5] = 59586
$data[22,0] = 45233
$data[22,1] = 1.06641
$data[22,2] = 1.06668
$data[22,3] = 1.06104
$data[22,4] = 1.06153
$data[22,5] = 74056
$data[23,0] = 45233
$data[23,1] = 1.06153
$data[23,2] = 1.06325
$data[23,3] = 1.06137
$data[23,4] = 1.06216
$data[23,5] = 27223
$data[24,0] = 45233
$data[24,1] = 1.06216
$data[24,2] = 1.06237
$data[24,3] = 1.06145
$data[24,4] = 1.06181
$data[24,5] = 9552
$data[25,0] = 45233
$data[25,1] = 1.06181
$data[25,2] = 1.06287
$data[25,3] = 1.06147
$data[25,4] = 1.06263
$data[25,5] = 14019
$data[26,0] = 45233
$data[26,1] = 1.06263
$data[26,2] = 1.0634
$data[26,3] = 1.06197
$data[26,4] = 1.0627
$data[26,5] = 23544
$data[27,0] = 45233
$data[27,1] = 1.0627
$data[27,2] = 1.07194
$data[27,3] = 1.06241
$data[27,4] = 1.07093
$data[27,5] = 59851
$data[28,0] = 45234
$data[28,1] = 1.07093
$data[28,2] = 1.07393
$data[28,3] = 1.06935
$data[28,4] = 1.07392
$data[28,5] = 106662
$data[29,0] = 45234
$data[29,1] = 1.07392
$data[29,2] = 1.07467
$data[29,3] = 1.07224
$data[29,4] = 1.07277
$data[29,5] = 26350
$data[30,0] = 45236
$data[30,1] = 1.07215
$data[30,2] = 1.07335
$data[30,3] = 1.07215
$data[30,4] = 1.07283
$data[30,5] = 13826
$data[31,0] = 45236
$data[31,1] = 1.07283
$data[31,2] = 1.07388
$data[31,3] = 1.07263
$data[31,4] = 1.07355
$data[31,5] = 14112
$data[32,0] = 45236
$data[32,1] = 1.07355
$data[32,2] = 1.07563
$data[32,3] = 1.07315
$data[32,4] = 1.07471
$data[32,5] = 35364
$data[33,0] = 45236
$data[33,1] = 1.07471
$data[33,2] = 1.07557
$data[33,3] = 1.0735
$data[33,4] = 1.07374
$data[33,5] = 34867
$data[34,0] = 45237
$data[34,1] = 1.07374
$data[34,2] = 1.07493
$data[34,3] = 1.0726
$data[34,4] = 1.07354
$data[34,5] = 43962
$data[35,0] = 45237
$data[35,1] = 1.07354
$data[35,2] = 1.0737
$data[35,3] = 1.07163
$data[35,4] = 1.07172
$data[35,5] = 22305
$data[36,0] = 45237
$data[36,1] = 1.07172
$data[36,2] = 1.07221
$data[36,3] = 1.07079
$data[36,4] = 1.07103
$data[36,5] = 12130
$data[37,0] = 45237
$data[37,1] = 1.07103
$data[37,2] = 1.0714
$data[37,3] = 1.07052
$data[37,4] = 1.07104
$data[37,5] = 16241
$data[38,0] = 45237
$data[38,1] = 1.07104
$data[38,2] = 1.07116
$data[38,3] = 1.0683
$data[38,4] = 1.06836
$data[38,5] = 34146
$data[39,0] = 45237
$data[39,1] = 1.06836
$data[39,2] = 1.06976
$data[39,3] = 1.06655
$data[39,4] = 1.06759
$data[39,5] = 41765
$data[40,0] = 45238
$data[40,1] = 1.06759
$data[40,2] = 1.06928
$data[40,3] = 1.0664
$data[40,4] = 1.06773
$data[40,5] = 52990
$data[41,0] = 45238
$data[41,1] = 1.06773
$data[41,2] = 1.07048
$data[41,3] = 1.0676
$data[41,4] = 1.06997
$data[41,5] = 24363
$data[42,0] = 45238
$data[42,1] = 1.06997
$data[42,2] = 1.06998
$data[42,3] = 1.06863
$data[42,4] = 1.06918
$data[42,5] = 13102
$data[43,0] = 45238
$data[43,1] = 1.06918
$data[43,2] = 1.06991
$data[43,3] = 1.06853
$data[43,4] = 1.06861
$data[43,5] = 15441
$data[44,0] = 45238
$data[44,1] = 1.06861
$data[44,2] = 1.0688
$data[44,3] = 1.06602
$data[44,4] = 1.06651
$data[44,5] = 33518
$data[45,0] = 45238
$data[45,1] = 1.06651
$data[45,2] = 1.06778
$data[45,3] = 1.06592
$data[45,4] = 1.06692
$data[45,5] = 35547
$data[46,0] = 45239
$data[46,1] = 1.06692
$data[46,2] = 1.0716
$data[46,3] = 1.06689
$data[46,4] = 1.07077
$data[46,5] = 55068
$data[47,0] = 45239
$data[47,1] = 1.07077
$data[47,2] = 1.07141
$data[47,3] = 1.06999
$data[47,4] = 1.0709
$data[47,5] = 26556
$data[48,0] = 45239
$data[48,1] = 1.0709
$data[48,2] = 1.07146
$data[48,3] = 1.07023
$data[48,4] = 1.0712
$data[48,5] = 12311
$data[49,0] = 45239
$data[49,1] = 1.0712
$data[49,2] = 1.07156
$data[49,3] = 1.07021
$data[49,4] = 1.07025
$data[49,5] = 13604
$data[50,0] = 45239
$data[50,1] = 1.07025
$data[50,2] = 1.07097
$data[50,3] = 1.06934
$data[50,4] = 1.06977
$data[50,5] = 31167
$data[51,0] = 45239
$data[51,1] = 1.06977
$data[51,2] = 1.07088
$data[51,3] = 1.068
$data[51,4] = 1.07075
$data[51,5] = 39342
$data[52,0] = 45240
$data[52,1] = 1.07075
$data[52,2] = 1.07255
$data[52,3] = 1.06998
$data[52,4] = 1.07047
$data[52,5] = 56595
$data[53,0] = 45240
$data[53,1] = 1.07047
$data[53,2] = 1.07056
$data[53,3] = 1.06601
$data[53,4] = 1.06673
$data[53,5] = 64566
$data[54,0] = 45240
$data[54,1] = 1.06673
$data[54,2] = 1.06704
$data[54,3] = 1.06641
$data[54,4] = 1.06653
$data[54,5] = 3460
$data[55,0] = 45240
$data[55,1] = 1.07075
$data[55,2] = 1.07255
$data[55,3] = 1.06998
$data[55,4] = 1.07047
$data[55,5] = 56603
$data[56,0] = 45240
$data[56,1] = 1.07047
$data[56,2] = 1.07056
$data[56,3] = 1.06601
$data[56,4] = 1.06673
$data[56,5] = 64726
$data[57,0] = 45240
$data[57,1] = 1.06673
$data[57,2] = 1.06751
$data[57,3] = 1.06629
$data[57,4] = 1.0672
$data[57,5] = 13983
$data[58,0] = 45240
$data[58,1] = 1.0672
$data[58,2] = 1.06744
$data[58,3] = 1.06657
$data[58,4] = 1.0669
$data[58,5] = 13820
$data[59,0] = 45240
$data[59,1] = 1.0669
$data[59,2] = 1.068
$data[59,3] = 1.06561
$data[59,4] = 1.06772
$data[59,5] = 35112
$data[60,0] = 45240
$data[60,1] = 1.06772
$data[60,2] = 1.06928
$data[60,3] = 1.06693
$data[60,4] = 1.06887
$data[60,5] = 38382
$data[61,0] = 45241
$data[61,1] = 1.06887
$data[61,2] = 1.06896
$data[61,3] = 1.06623
$data[61,4] = 1.06792
$data[61,5] = 56460
$data[62,0] = 45241
$data[62,1] = 1.06792
$data[62,2] = 1.06875
$data[62,3] = 1.06729
$data[62,4] = 1.06842
$data[62,5] = 20019
$data[63,0] = 45243
$data[63,1] = 1.0683
$data[63,2] = 1.06907
$data[63,3] = 1.06805
$data[63,4] = 1.06867
$data[63,5] = 13886
$data[64,0] = 45243
$data[64,1] = 1.06867
$data[64,2] = 1.06873
$data[64,3] = 1.06814
$data[64,4] = 1.06855
$data[64,5] = 6473

# Copy the existing date-column formatting (style used by column A cells,
# e.g. the datetime number format) down into the new rows before writing
# the values, so the new "datetime" column cells match the rest of the table.
$ws.Range("A4992").Copy() | Out-Null
$ws.Range("A4993:A5057").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Bulk-write the new rows in one shot.
$ws.Range("A4993:F5057").Value = $data
